# actualizacion de lista de precios
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B12" = 4800.0
    "B13" = 4390.0
    "B14" = 4260.0
    "B15" = 3755.0
    "B16" = 4850.0
    "B23" = 8580.0
    "B56" = 5480.0
    "B57" = 5010.0
    "B58" = 8390.0
    "B59" = 7310.0
    "B60" = 8735.0
    "B61" = 8735.0
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
